# Add four new dtype/ECV vocabulary mapping rows to Sheet1, continuing
# directly after the existing last row (85).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("dtype uncert", "broader", "cciecv_landCov", "narrower"),
    @("dtype_vp", "broader", "cciecv_vegParam", "narrower"),
    @("dtype_wl", "broader", "cciecv_riverDischarge", "narrower"),
    @("dtype_lcchange", "broader", "cciecv_landCov", "narrower")
)

$startRow = 86
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    $ws.Range("A$r").Value = $values[0]
    $ws.Range("B$r").Value = $values[1]
    $ws.Range("C$r").Value = $values[2]
    $ws.Range("D$r").Value = $values[3]
}
